$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 35-47 get their Id (A), Antal (I), Ost (Q) and Nord (R) values
# re-assigned; every other column in those rows is left untouched.
# (The underlying data got re-ordered/re-matched across these rows -
# the same four fields simply end up attached to a different row.)
$rows = @(
    @{ Row = 35; A = 111675572; I = "2"; Q = 690494.5947179901; R = 6661104.692649405 },
    @{ Row = 36; A = 111675575; I = "1"; Q = 690480.7418955797; R = 6661091.463633558 },
    @{ Row = 37; A = 111675579; I = "3"; Q = 690425.8424831247; R = 6661357.862056008 },
    @{ Row = 38; A = 111675573; I = "2"; Q = 690487.9917822112; R = 6661106.352564453 },
    @{ Row = 39; A = 111675582; I = "1"; Q = 690352.3333891984; R = 6661470.655078794 },
    @{ Row = 40; A = 111675577; I = "1"; Q = 690430.9193086301; R = 6661356.623615522 },
    @{ Row = 41; A = 111675583; I = "1"; Q = 690415.8809986882; R = 6661424.403280765 },
    @{ Row = 42; A = 111675578; I = "1"; Q = 690368.3990222017; R = 6661295.837351476 },
    @{ Row = 43; A = 111675581; I = "1"; Q = 690413.7262835158; R = 6661427.29424896  },
    @{ Row = 44; A = 111675584; I = "2"; Q = 690414.984509701;  R = 6661422.355185229 },
    @{ Row = 45; A = 111675574; I = "1"; Q = 690486.6986671695; R = 6661102.281881573 },
    @{ Row = 46; A = 111675580; I = "3"; Q = 690370.5537696742; R = 6661292.946251329 },
    @{ Row = 47; A = 111675571; I = "5"; Q = 690509.4285896254; R = 6661040.900344189 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A

    # Column I ("Antal") is stored as text in this sheet. Force the cell to
    # text before writing so Excel doesn't auto-convert the numeric-looking
    # string into a Number, then drop the temporary number format again so
    # no stray style/quote-prefix is left behind on the cell.
    $iCell = $ws.Cells.Item($r.Row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = $r.I
    $iCell.Style = "Normal"

    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
}
